# Update the id values in column A (rows 2 and 3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 6
$ws.Range("A3").Value = 8

# Move / record the current selection to C4, matching the saved view state
$ws.Range("C4").Select()
